$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vessels")

# Delete row 2 (duplicate header row), shifting all data rows up by one
$ws.Rows.Item(2).Delete()

# Restore the active selection to match target state
$ws.Range("K9").Select()
